# NN needs to be corrected
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @{ B = 1234.031127929688; C = 0.9258; D = 0.9351000189781189; E = 1.250499963760376; F = 0.5346999764442444; H = 0.8289 }
    3 = @{ B = 1193.449829101562; C = 0.9548; D = 0.9539;              E = 1.06659996509552;  F = 0.7386000156402588; H = 0.9953 }
    4 = @{ B = 815.7728881835938; C = 0.953;  D = 0.9497;              E = 1.119099974632263; F = 0.756600022315979;  H = 0.9586 }
    5 = @{ B = 810.4166259765625; C = 0.8549; D = 0.8588;              E = 1.004799962043762; F = 0.449999988079071;  H = 0.1534 }
    6 = @{ B = 1119.871948242188; C = 0.886;  D = 0.8921;              E = 1.035400032997131; F = 0.5665000081062317; H = 0.4484 }
    7 = @{ B = 874.6599731445312; C = 0.8808; D = 0.8798999786376953;  E = 1.01830005645752;  F = 0.7138000130653381; H = 0.3398 }
    8 = @{ B = 970.0062866210938; C = 0.8692; D = 0.8683;              E = 1.023699998855591; F = 0.7228000164031982; H = 0.2375 }
    9 = @{ B = 7018.208984375;    C = 0.9044; D = 0.9121;              E = 1.250499963760376; F = 0.449999988079071;  H = 3.9619 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("B$row").Value = $vals.B
    $ws.Range("C$row").Value = $vals.C
    $ws.Range("D$row").Value = $vals.D
    $ws.Range("E$row").Value = $vals.E
    $ws.Range("F$row").Value = $vals.F
    $ws.Range("H$row").Value = $vals.H
}
